# Apply the "Journal de travail" update:
#  - log two new work entries (27/06/2023 and 28/06/2023)
#  - update the sheet's active selection / scroll position

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")
$ws.Activate()

# Row 61 : 27/06/2023 - Implémentation - 6h - Frontend: Refactor et fonctionnalités de personnalisation
$ws.Range("A61").Value = 45104
$ws.Range("B61").Value = "Implémentation"
$ws.Range("C61").Value = 6
$ws.Range("D61").Value = "Frontend: Refactor et fonctionnalités de personnalisation"

# Row 62 : 28/06/2023 - Implémentation - Frontend: Correctifs
$ws.Range("A62").Value = 45105
$ws.Range("B62").Value = "Implémentation"
$ws.Range("D62").Value = "Frontend: Correctifs"

# Update view: scroll so column A is the left-most visible column, and
# move the active selection to D65.
$win = $excel.ActiveWindow
$win.ScrollRow = 44
$win.ScrollColumn = 1
$ws.Range("D65").Select()
